# Rerun analyses for misinfotext: add a "2015" sheet (inserted between
# "2014" and "2016") with its own factcheck/polarity/subjectivity data +
# hyperlinks, and update the "Summary" sheet with a new 2015 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2015" worksheet right before "2016".
# ---------------------------------------------------------------------
$ws2016 = $wb.Worksheets.Item("2016")
$ws2015 = $wb.Worksheets.Add($ws2016)
$ws2015.Name = "2015"

# Header row (reuses the same shared strings as the other year sheets).
$ws2015.Range("A1").Value = "factcheckURL"
$ws2015.Range("B1").Value = "polarity"
$ws2015.Range("C1").Value = "subjectivity"
$ws2015.Range("A1:C1").Style = "Normal"

# Data rows: factcheckURL, polarity, subjectivity.
$data2015 = @(
    @("https://www.politifact.com/factchecks/2015/nov/06/greg-abbott/greg-abbott-embarrassed-says-californians-buying-m/", -0.1318181818181818, 0.5772727272727273),
    @("https://www.politifact.com/factchecks/2015/oct/16/scott-walker/gov-scott-walker-says-board-wanted-accept-mickey-m/", 0, 0),
    @("https://www.politifact.com/factchecks/2015/oct/05/dana-loesch/Planned-parenthood-86-percent-abortion-revenue/", -0.1, 0.6),
    @("https://www.politifact.com/factchecks/2015/sep/23/donald-trump/hillary-clinton-obama-birther-fact-check/", 0, 0),
    @("https://www.politifact.com/factchecks/2015/sep/10/ted-cruz/ted-cruz-says-deal-will-facilitate-and-accelerate-/", 0, 0),
    @("https://www.politifact.com/factchecks/2015/aug/26/hillary-clinton/hillary-clinton-says-no-gop-candidate-has-talked-a/", -0.475, 0.6),
    @("https://www.politifact.com/factchecks/2015/jul/21/wisconsin-state-afl-cio/wisconsin-afl-cio-says-scott-walker-budget-means-n/", 0, 0),
    @("https://www.politifact.com/factchecks/2015/mar/20/glenn-beck/glenn-beck-says-barack-obama-took-iran-hamas-us-te/", 0.07812424843674844, 0.4032627865961199),
    @("https://www.politifact.com/factchecks/2015/nov/08/ben-carson/ben-carson-said-no-one-who-signed-declaration-inde/", 0.07460724759111856, 0.4004198668714798),
    @("https://www.politifact.com/factchecks/2015/jul/07/dinesh-dsouza/hillary-clinton-confederate-battle-flag-nope-old-i/", -0.4, 0.7),
    @("https://www.politifact.com/factchecks/2015/jun/25/gavin-mcinnes/tweet-civil-war-was-about-secession-not-slavery/", 0.05, 0.1),
    @("https://www.politifact.com/factchecks/2015/feb/06/scott-walker/despite-deliberate-actions-scott-walker-calls-chan/", 0, 1)
)

for ($i = 0; $i -lt $data2015.Length; $i++) {
    $row = $i + 2
    $entry = $data2015[$i]
    $urlCell = $ws2015.Cells.Item($row, 1)
    $urlCell.Value = $entry[0]
    $ws2015.Hyperlinks.Add($urlCell, $entry[0])
    $urlCell.Style = "Hyperlink"
    $ws2015.Cells.Item($row, 2).Value = $entry[1]
    $ws2015.Cells.Item($row, 3).Value = $entry[2]
}

# ---------------------------------------------------------------------
# 2. Update the "Summary" sheet: insert a 2015 row between 2014 and 2016.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Rows.Item(6).Insert()

# Copy formatting from the row above (2014) so the new row matches style.
$summary.Range("A5:K5").Copy()
$summary.Range("A6:K6").PasteSpecial(-4122)

$summary.Range("A6").Value = 2015
$summary.Range("B6").Value = 0.07812424843674844
$summary.Range("C6").Value = -0.475
$summary.Range("D6").Value = -0.0753405571491929
$summary.Range("E6").Value = 0.1091291398205041
$summary.Range("F6").Value = 1
$summary.Range("G6").Value = 0
$summary.Range("H6").Value = 0.3650796150616939
$summary.Range("I6").Value = 3
$summary.Range("J6").Value = 4
$summary.Range("K6").Value = 5

Write-Host "2015 sheet inserted and Summary updated"
